$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily records (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti)
$data = @(
    @(44313, 7, 21, 186.0217911241031),
    @(44314, 0, 21, 186.0217911241031),
    @(44315, 2, 21, 186.0217911241031),
    @(44316, 2, 21, 186.0217911241031),
    @(44317, 4, 18, 159.4472495349455),
    @(44318, 3, 20, 177.1636105943839)
)

$startRow = 239

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]

    # Copy formatting from the column-A cell directly above (style index "2":
    # centered date format, border, bold) onto the new date cell.
    $ws.Range("A" + ($row - 1)).Copy()
    $ws.Range("A" + $row).PasteSpecial(-4122)  # xlPasteFormats

    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
}

$excel.CutCopyMode = 0
